# Apply updated crypto price/volume figures to Sheet1 (D = Price, E = Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.264.46"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "2.522.64"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'323.38"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'109.13"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.562"
$ws.Range("E9").Value = "  +4.41%  "
$ws.Range("D10").Value = "'40.43"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").Value = "'20.13"
$ws.Range("E11").Value = "  +8.73%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "'7.26"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "2.917.04"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "2.539.16"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "48.151.96"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "'13.21"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "0.0₃0946"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").Value = "'72.53"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "'268.50"
$ws.Range("E24").Value = "  +8.02%  "
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'26.22"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D29").Value = "'10.17"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "'0.146"
$ws.Range("E30").Value = "  +5.01%  "
$ws.Range("D31").Value = "'34.98"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "'20.03"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'5.40"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "'22.19"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "'118.83"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D45").Value = "2.002.01"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "'3.13"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  +6.87%  "
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "'9.08"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'80.47"
$ws.Range("E51").Value = "  +3.00%  "
